$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.63"
$ws.Range("E2").Value = "'1.60%"
$ws.Range("D3").Value = "'35.55"
$ws.Range("E3").Value = "'-1.97%"
$ws.Range("D4").Value = "'5.116"
$ws.Range("E4").Value = "'1.45%"
$ws.Range("D5").Value = "'0.08214"
$ws.Range("E5").Value = "'4.30%"
$ws.Range("D6").Value = "'2.071"
$ws.Range("E6").Value = "'-8.85%"
$ws.Range("D7").Value = "'7.993"
$ws.Range("E7").Value = "'-0.07%"
$ws.Range("D8").Value = "'4.130"
$ws.Range("E8").Value = "'-0.33%"
$ws.Range("D9").Value = "'2.981"
$ws.Range("E9").Value = "'12.07%"
$ws.Range("D10").Value = "'0.9280"
$ws.Range("E10").Value = "'0.03%"
$ws.Range("D11").Value = "'0.1080"
$ws.Range("E11").Value = "'10.42%"
$ws.Range("D12").Value = "'0.1924"
$ws.Range("E12").Value = "'3.25%"
$ws.Range("D13").Value = "'0.09682"
$ws.Range("E13").Value = "'7.90%"
$ws.Range("D14").Value = "'0.03639"
$ws.Range("E14").Value = "'-3.27%"
$ws.Range("D15").Value = "'0.09935"
$ws.Range("E15").Value = "'0.31%"
$ws.Range("D16").Value = "'0.001433"
$ws.Range("E16").Value = "'-0.24%"
$ws.Range("D17").Value = "'0.005676"
$ws.Range("E17").Value = "'0.68%"
$ws.Range("D18").Value = "'3.471"
$ws.Range("E18").Value = "'0.31%"
$ws.Range("D19").Value = "'0.3426"
$ws.Range("E19").Value = "'1.76%"
$ws.Range("D20").Value = "'0.1311"
$ws.Range("E20").Value = "'-0.64%"
$ws.Range("D21").Value = "'5.097"
$ws.Range("E21").Value = "'-0.38%"
$ws.Range("E22").Value = "'-2.05%"
$ws.Range("D23").Value = "'0.04542"
$ws.Range("E23").Value = "'-0.71%"
$ws.Range("D24").Value = "'0.001227"
$ws.Range("E24").Value = "'-0.65%"
$ws.Range("D25").Value = "'0.004808"
$ws.Range("E25").Value = "'0.77%"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'-3.86%"
$ws.Range("D27").Value = "'0.0004454"
$ws.Range("E27").Value = "'-6.01%"
$ws.Range("D39").Value = "'0.01995"
$ws.Range("E39").Value = "'4.16%"
$ws.Range("D40").Value = "'0.04936"
$ws.Range("E40").Value = "'-1.35%"
$ws.Range("D41").Value = "'0.007667"
$ws.Range("E41").Value = "'-1.82%"
$ws.Range("D42").Value = "'0.009841"
$ws.Range("E42").Value = "'25.84%"
$ws.Range("D43").Value = "'0.1387"
$ws.Range("E43").Value = "'-0.13%"
$ws.Range("D44").Value = "'0.002117"
$ws.Range("E44").Value = "'-1.20%"
$ws.Range("D45").Value = "'0.01160"
$ws.Range("E45").Value = "'1.65%"
$ws.Range("D46").Value = "'0.00006553"
$ws.Range("E46").Value = "'6.67%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.02%"
$ws.Range("D48").Value = "'61.91"
$ws.Range("E48").Value = "'19.62%"
$ws.Range("E49").Value = "'-16.84%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.02%"
